$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New attendance-log rows uploaded by the attendance app.
# Columns: Student ID | Subject | Log Date | Log Time | Type | User
# ---------------------------------------------------------------------
$records = @(
    @{ Row = 2;  Id = "201218"; Fill = $true  },
    @{ Row = 3;  Id = "211169"; Fill = $false },
    @{ Row = 4;  Id = "200943"; Fill = $true  },
    @{ Row = 5;  Id = "211121"; Fill = $false },
    @{ Row = 6;  Id = "191055"; Fill = $true  },
    @{ Row = 7;  Id = "202022"; Fill = $false },
    @{ Row = 8;  Id = "201513"; Fill = $true  },
    @{ Row = 9;  Id = "200997"; Fill = $false },
    @{ Row = 10; Id = "201880"; Fill = $true  }
)

$subject = "general surgery"
$logDate = "19/10/2025"
$logTime = "10:30:00"
$logType = "Excuse"
$logUser = "System"

foreach ($rec in $records) {
    $r = $rec.Row

    # Write the Student ID as TEXT (not a number) via a scratch cell +
    # TEXT() formula, then paste the resulting value only - this avoids
    # Excel auto-coercing the numeric-looking id into a Number cell.
    $ws.Range("Z1").Formula = '=TEXT(' + $rec.Id + ',"0")'
    $ws.Range("Z1").Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4163)
    $ws.Range("Z1").Clear()

    $ws.Cells.Item($r, 2).Value = $subject
    $ws.Cells.Item($r, 3).Value = $logDate
    $ws.Cells.Item($r, 4).Value = $logTime
    $ws.Cells.Item($r, 5).Value = $logType
    $ws.Cells.Item($r, 6).Value = $logUser
}

# ---------------------------------------------------------------------
# Row banding: even data rows keep the existing shaded style (copied
# from row 2), odd data rows use the same font/alignment without the
# background fill.
# ---------------------------------------------------------------------
$ws.Range("A2:F2").Copy()
foreach ($rec in $records) {
    if ($rec.Row -ne 2) {
        $target = $ws.Range("A" + $rec.Row + ":F" + $rec.Row)
        $target.PasteSpecial(-4122)
        if (-not $rec.Fill) {
            $target.Interior.Pattern = -4142
        }
    }
}

$excel.CutCopyMode = $false
